$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# age (B), sex (C), edu (D) values for subjects 1-50 (rows 2-51)
$data = @(
    @(53, "Male", 5),
    @(51, "Male", 7),
    @(78, "Female", 7),
    @(58, "Male", 7),
    @(72, "Male", 12),
    @(58, "Male", 10),
    @(57, "Female", 9),
    @(69, "Male", 7),
    @(50, "Female", 12),
    @(56, "Male", 7),
    @(74, "Female", 10),
    @(77, "Male", 5),
    @(79, "Male", 6),
    @(73, "Male", 7),
    @(52, "Male", 8),
    @(63, "Male", 10),
    @(66, "Male", 7),
    @(69, "Male", 7),
    @(63, "Male", 5),
    @(72, "Female", 6),
    @(62, "Male", 15),
    @(78, "Male", 13),
    @(60, "Male", 9),
    @(72, "Male", 7),
    @(58, "Female", 13),
    @(56, "Female", 7),
    @(60, "Male", 9),
    @(68, "Female", 5),
    @(67, "Female", 5),
    @(57, "Male", 7),
    @(77, "Female", 11),
    @(50, "Male", 6),
    @(66, "Female", 7),
    @(58, "Male", 10),
    @(50, "Male", 13),
    @(76, "Female", 8),
    @(77, "Male", 7),
    @(66, "Male", 13),
    @(71, "Male", 10),
    @(56, "Male", 7),
    @(50, "Female", 6),
    @(77, "Male", 6),
    @(70, "Female", 10),
    @(73, "Female", 11),
    @(69, "Female", 5),
    @(54, "Male", 8),
    @(76, "Female", 5),
    @(53, "Male", 5),
    @(74, "Female", 12),
    @(79, "Male", 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

